$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 48/49 swap the Cronos and BinanceUSD entries; every other row keeps its
# coin/link and only refreshes Price (D) and/or Volume(1h) (E) with the latest scrape.

$ws.Range('D2').Value = '40.825.38'
$ws.Range('E2').Value = '  -1.88%  '

$ws.Range('D3').Value = '2.163.81'
$ws.Range('E3').Value = '  -2.86%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.07%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '246.57'
$ws.Range('E5').Value = '  -2.84%  '

$ws.Range('E6').Value = '  -2.29%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '66.06'
$ws.Range('E7').Value = '  -6.40%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  -0.02%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.565'
$ws.Range('E9').Value = '  +0.53%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '58.13'
$ws.Range('E10').Value = '  -0.71%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '35.65'
$ws.Range('E11').Value = '  -15.09%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0922'
$ws.Range('E12').Value = '  -4.38%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.103'
$ws.Range('E13').Value = '  -1.32%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.87'
$ws.Range('E14').Value = '  -0.81%  '

$ws.Range('D15').Value = '2.483.62'
$ws.Range('E15').Value = '  -2.84%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.857'
$ws.Range('E16').Value = '  +0.22%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.24'
$ws.Range('E17').Value = '  -4.64%  '

$ws.Range('D18').Value = '2.166.26'
$ws.Range('E18').Value = '  -3.01%  '

$ws.Range('D19').Value = '40.753.79'
$ws.Range('E19').Value = '  -1.94%  '

$ws.Range('D20').Value = '0.0₃0936'
$ws.Range('E20').Value = '  -2.99%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.07'
$ws.Range('E21').Value = '  -1.76%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '71.15'
$ws.Range('E22').Value = '  -2.34%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '228.17'
$ws.Range('E23').Value = '  -2.66%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.13'
$ws.Range('E24').Value = '  -6.04%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.65'
$ws.Range('E25').Value = '  +14.04%  '

$ws.Range('E26').Value = '  +0.11%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.70'
$ws.Range('E27').Value = '  -1.40%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.42'
$ws.Range('E28').Value = '  -4.11%  '

$ws.Range('E29').Value = '  -5.90%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '168.58'
$ws.Range('E30').Value = '  -1.52%  '

$ws.Range('E31').Value = '  -5.18%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.11'
$ws.Range('E32').Value = '  -2.42%  '

$ws.Range('E33').Value = '  -1.11%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.56'
$ws.Range('E34').Value = '  +1.24%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0736'
$ws.Range('E35').Value = '  +2.32%  '

$ws.Range('E36').Value = '  -3.30%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.52'
$ws.Range('E37').Value = '  -3.54%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '24.67'
$ws.Range('E38').Value = '  -7.02%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.94'
$ws.Range('E39').Value = '  -2.52%  '

$ws.Range('E40').Value = '  +4.10%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.17'
$ws.Range('E41').Value = '  -5.16%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.43'
$ws.Range('E42').Value = '  -9.47%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '11.57'
$ws.Range('E43').Value = '  -3.41%  '

$ws.Range('E44').Value = '  -6.00%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '60.19'
$ws.Range('E45').Value = '  -14.07%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.191'
$ws.Range('E46').Value = '  -9.22%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.46'
$ws.Range('E47').Value = '  -3.78%  '

$ws.Range('E50').Value = '  -0.87%  '

$ws.Range('E51').Value = '  -3.74%  '

# Rows 48 and 49 swap places (Cronos <-> BinanceUSD), each with a freshly updated Volume(1h)
$ws.Range('B48').Value = 'BinanceUSD'
$ws.Range('C48').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.00'
$ws.Range('E48').Value = '  -0.16%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0989'
$ws.Range('E49').Value = '  -2.75%  '
